$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "246.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "0.58%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "26.40"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "5.72%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.090"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.84%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05603"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.37%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.478"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-0.94%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.8135"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.35%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8454"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.81%"
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1336"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.02%"
$ws.Range("B10").Value = "BitrueCoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.02857"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.89%"
$ws.Range("B11").Value = "BitMartToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09389"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-0.19%"
$ws.Range("B12").Value = "BitForexToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.001519"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.12%"
$ws.Range("B13").Value = "One"
$ws.Range("C13").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0005961"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-93.87%"
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.006159"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.81%"
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.595"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.75%"
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.011"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.93%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.055"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.09%"
$ws.Range("B18").Value = "BitpandaEcosystemToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3207"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.65%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07005"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.99%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03178"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-2.86%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.49%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.750"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04654"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.87%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001246"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.24%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004583"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "1.19%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.91%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "168.15%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03668"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.11%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1057"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "0.51%"
$ws.Range("B42").Value = "KickToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.006211"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.20%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002500"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-8.08%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.008794"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "5.38%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005297"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.24%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.12%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-42.03%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002648"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "29.37%"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.12%"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.12%"
